$wb = $excel.ActiveWorkbook

# StatOutput sheet: update the 4 numeric-looking result cells in row 2.
# Values that look like numbers get auto-coerced to the Number type by
# plain Value assignment, so stage them as text via a formula in a scratch
# cell and paste-special (values) into place to keep them as shared strings.
$statOutput = $wb.Worksheets.Item("StatOutput")

$statOutput.Range("Z1").Formula = '="2"'
$statOutput.Range("Z1").Copy()
$statOutput.Range("A2").PasteSpecial()

$statOutput.Range("Z1").Formula = '="5"'
$statOutput.Range("Z1").Copy()
$statOutput.Range("B2").PasteSpecial()

$statOutput.Range("Z1").Formula = '="2"'
$statOutput.Range("Z1").Copy()
$statOutput.Range("C2").PasteSpecial()

$statOutput.Range("Z1").Value = ""

# StatOutput_Message sheet: update the Cypher query text on row 18 (Akita -> Greyhound)
$statOutputMessage = $wb.Worksheets.Item("StatOutput_Message")
$statOutputMessage.Range("A18").Value = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Greyhound']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"
